$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 377.33334
$ws.Range("I2").Value = 252.8
$ws.Range("K2").Value = 252.8
$ws.Range("M2").Value = -139.8

$ws.Range("H33").Value = 292.875
$ws.Range("I33").Value = 263.83334
$ws.Range("K33").Value = 263.83334
$ws.Range("M33").Value = -34.83334000000002

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8709.799999999999
$ws.Range("I32").Value = 6933.321
$ws.Range("J32").Value = 22160.285
$ws.Range("K32").Value = 6933.321
$ws.Range("L32").Value = 22160.285
$ws.Range("M32").Value = -6646.321
$ws.Range("N32").Value = -22734.285

$ws.Range("H74").Value = 103029.055
$ws.Range("I74").Value = 123913.93
$ws.Range("J74").Value = 16506
$ws.Range("K74").Value = 123913.93
$ws.Range("L74").Value = 16506
$ws.Range("M74").Value = -123039.93
$ws.Range("N74").Value = -18254

$ws.Range("H77").Value = 103029.055
$ws.Range("I77").Value = 123913.93
$ws.Range("J77").Value = 16506
$ws.Range("K77").Value = 619569.6499999999
$ws.Range("L77").Value = 82530
$ws.Range("M77").Value = -615201.6499999999
$ws.Range("N77").Value = -91266

$ws.Range("H132").Value = 5436.7896
$ws.Range("I132").Value = 1466.619
$ws.Range("J132").Value = 10341.117
$ws.Range("K132").Value = 4399.857
$ws.Range("L132").Value = 31023.351
$ws.Range("M132").Value = -1869.857
$ws.Range("N132").Value = -36083.351

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2907.9783
$ws.Range("I31").Value = 2117.1482
$ws.Range("J31").Value = 4031.7896
$ws.Range("K31").Value = 2117.1482
$ws.Range("L31").Value = 4031.7896
$ws.Range("M31").Value = -1822.1482
$ws.Range("N31").Value = -4621.7896

$ws.Range("H34").Value = 2907.9783
$ws.Range("I34").Value = 2117.1482
$ws.Range("J34").Value = 4031.7896
$ws.Range("K34").Value = 2117.1482
$ws.Range("L34").Value = 4031.7896
$ws.Range("M34").Value = -1915.1482
$ws.Range("N34").Value = -4435.7896

$ws.Range("H58").Value = 2117042.8
$ws.Range("I58").Value = 3136894.2
$ws.Range("J58").Value = 4493.357
$ws.Range("K58").Value = 3136894.2
$ws.Range("L58").Value = 4493.357
$ws.Range("M58").Value = -3136691.2
$ws.Range("N58").Value = -4899.357

$ws.Range("H132").Value = 2470.4849
$ws.Range("I132").Value = 1855.7391
$ws.Range("J132").Value = 3884.4
$ws.Range("K132").Value = 5567.2173
$ws.Range("L132").Value = 11653.2
$ws.Range("M132").Value = -3037.2173
$ws.Range("N132").Value = -16713.2

$ws.Range("H136").Value = 2117042.8
$ws.Range("I136").Value = 3136894.2
$ws.Range("J136").Value = 4493.357
$ws.Range("K136").Value = 9410682.600000001
$ws.Range("L136").Value = 13480.071
$ws.Range("M136").Value = -9408132.600000001
$ws.Range("N136").Value = -18580.071

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 9710.454
$ws.Range("I68").Value = 550.1667
$ws.Range("J68").Value = 20702.8
$ws.Range("K68").Value = 1650.5001
$ws.Range("L68").Value = 62108.39999999999
$ws.Range("M68").Value = -839.5001
$ws.Range("N68").Value = -63730.39999999999

$ws.Range("H71").Value = 9710.454
$ws.Range("I71").Value = 550.1667
$ws.Range("J71").Value = 20702.8
$ws.Range("K71").Value = 4951.5003
$ws.Range("L71").Value = 186325.2
$ws.Range("M71").Value = -895.5002999999997
$ws.Range("N71").Value = -194437.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 15278.7
$ws.Range("I46").Value = 8020.5
$ws.Range("J46").Value = 17093.25
$ws.Range("K46").Value = 8020.5
$ws.Range("L46").Value = 17093.25
$ws.Range("M46").Value = -7864.5
$ws.Range("N46").Value = -17405.25

$ws.Range("H132").Value = 10714.417
$ws.Range("I132").Value = 18526.834
$ws.Range("J132").Value = 2902
$ws.Range("K132").Value = 55580.50199999999
$ws.Range("L132").Value = 8706
$ws.Range("M132").Value = -53050.50199999999
$ws.Range("N132").Value = -13766

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 4900
$ws.Range("I53").Value = 4900
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 4900
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -4382
$ws.Range("N53").ClearContents()

$ws.Range("H63").Value = 34014.168
$ws.Range("J63").Value = 34014.168
$ws.Range("L63").Value = 34014.168
$ws.Range("N63").Value = -35512.168

$ws.Range("H66").Value = 34014.168
$ws.Range("J66").Value = 34014.168
$ws.Range("L66").Value = 102042.504
$ws.Range("N66").Value = -109530.504

$ws.Range("H70").Value = 30097.25
$ws.Range("J70").Value = 30097.25
$ws.Range("L70").Value = 30097.25
$ws.Range("N70").Value = -30637.25

$ws.Range("H73").Value = 30097.25
$ws.Range("J73").Value = 30097.25
$ws.Range("L73").Value = 30097.25
$ws.Range("N73").Value = -31969.25

$ws.Range("H136").Value = 6372.1875
$ws.Range("I136").Value = 4756.875
$ws.Range("J136").Value = 7987.5
$ws.Range("K136").Value = 14270.625
$ws.Range("L136").Value = 23962.5
$ws.Range("M136").Value = -11720.625
$ws.Range("N136").Value = -29062.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 11500
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 11500
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 11500
$ws.Range("N26").Value = -12086
$ws.Range("M26").ClearContents()

$ws.Range("H29").Value = 10600
$ws.Range("I29").Value = 8600
$ws.Range("J29").Value = 11000
$ws.Range("K29").Value = 8600
$ws.Range("L29").Value = 11000
$ws.Range("M29").Value = -8310
$ws.Range("N29").Value = -11580

$ws.Range("H30").Value = 9560
$ws.Range("J30").Value = 8500
$ws.Range("L30").Value = 8500
$ws.Range("N30").Value = -8714

$ws.Range("H31").Value = 10100
$ws.Range("I31").Value = 7500
$ws.Range("J31").Value = 10750
$ws.Range("K31").Value = 7500
$ws.Range("L31").Value = 10750
$ws.Range("M31").Value = -7152
$ws.Range("N31").Value = -11446

$ws.Range("H126").Value = 1209.4546
$ws.Range("I126").Value = 1033.7778
$ws.Range("K126").Value = 3101.3334
$ws.Range("M126").Value = -631.3334000000004

$ws.Range("H136").Value = 4766.017
$ws.Range("I136").Value = 2109.5312
$ws.Range("J136").Value = 7914.4443
$ws.Range("K136").Value = 6328.5936
$ws.Range("L136").Value = 23743.3329
$ws.Range("M136").Value = -3778.5936
$ws.Range("N136").Value = -28843.3329
